$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff Datetime / Correspond Handback DateTime for row 3
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E3").Value = "2016-03-21 03:14:44"
$wsZh.Range("H3").Value = "2016-03-21 03:15:29"

# de-de sheet: update Correspond Handoff Datetime / Correspond Handback DateTime for row 3
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E3").Value = "2016-03-21 03:14:52"
$wsDe.Range("H3").Value = "2016-03-21 03:15:43"
